# LOE updates for August
# Enter the August hours (column H) for each team member's worksheet.
# Each worksheet has a weekly/monthly LOE grid; column H corresponds to the
# 08/15/09 period. After entry, all dependent formulas (row 18 totals,
# row 21/22/23 remaining-LOE calcs, the Project Summary / Milestones roll-ups)
# recalculate automatically.

$wb = $excel.ActiveWorkbook

# Each entry is (SheetName, Row, Value) for column H (the 08/15/09 / August
# period). SB-Dev and SB-Manager are applied first: their row-21/22 "days
# remaining" formulas read 'Project Summary'!H16 (a roll-up of every other
# sheet's H18) inside an IF(), and this keeps that recalculation well-formed
# once every sheet's H column has been entered.
$order = @(
    @("SB-Dev", 11, 160),
    @("SB-Manager", 6, 16),
    @("Warren", 6, 12), @("Warren", 7, 8), @("Warren", 14, 8),
    @("Sean", 6, 50), @("Sean", 7, 10), @("Sean", 14, 10), @("Sean", 16, 10),
    @("Rhett", 10, 128),
    @("Lee", 7, 4), @("Lee", 13, 4),
    @("Jalpa", 10, 80), @("Jalpa", 12, 80),
    @("David", 7, 20), @("David", 13, 20),
    @("Jignesh", 15, 16), @("Jignesh", 17, 16),
    @("Nataliya", 8, 80), @("Nataliya", 13, 80),
    @("John", 10, 160),
    @("Renee", 7, 8),
    @("Dong", 7, 12), @("Dong", 14, 4), @("Dong", 15, 6), @("Dong", 16, 6), @("Dong", 17, 5)
)

foreach ($item in $order) {
    $ws = $wb.Worksheets.Item($item[0])
    $ws.Cells.Item($item[1], 8).Value = $item[2]
}

# Restore the cursor/selection on each sheet to match where the editor last
# left off (cosmetic, matches the final authored state).
$selections = @{
    "Warren"     = "H12"
    "Sean"       = "H8"
    "Rhett"      = "I10"
    "Lee"        = "H8"
    "Jalpa"      = "H11"
    "David"      = "I7"
    "Jignesh"    = "I15"
    "Nataliya"   = "I8"
    "John"       = "H11"
    "Renee"      = "H8"
    "Dong"       = "H24"
    "SB-Dev"     = "H12"
    "SB-Manager" = "H7"
}

foreach ($sheetName in $selections.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Activate()
    $ws.Range($selections[$sheetName]).Select()
}

# SB-Manager was the tab on-screen when the file was last saved; make sure it
# ends up active again (and scrolled so row 6 is at the top, per the source).
$sbManager = $wb.Worksheets.Item("SB-Manager")
$sbManager.Activate()
$excel.ActiveWindow.ScrollRow = 6
$sbManager.Range("H7").Select()
